$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3617.158
$ws.Range("I17").Value = 1399
$ws.Range("J17").Value = 3740.389
$ws.Range("K17").Value = 4197
$ws.Range("L17").Value = 11221.167
$ws.Range("M17").Value = -4029
$ws.Range("N17").Value = -11557.167
$ws.Range("H28").Value = 478.5
$ws.Range("I28").Value = 313.41177
$ws.Range("K28").Value = 313.41177
$ws.Range("M28").Value = 171.58823
$ws.Range("H29").Value = 9096.200000000001
$ws.Range("J29").Value = 10245.5
$ws.Range("L29").Value = 30736.5
$ws.Range("N29").Value = -31298.5
$ws.Range("H33").Value = 799
$ws.Range("I33").Value = 376.55554
$ws.Range("K33").Value = 376.55554
$ws.Range("M33").Value = -147.55554
$ws.Range("H38").Value = 594.0909
$ws.Range("I38").Value = 594.0909
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1782.2727
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -1410.2727
$ws.Range("H40").Value = 2536.0908
$ws.Range("I40").Value = 2098.3333
$ws.Range("K40").Value = 2098.3333
$ws.Range("M40").Value = -1923.3333
$ws.Range("H87").Value = 174999.5
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 174999.5
$ws.Range("K87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("M87").Value = 174999.5
$ws.Range("N87").Value = -177495.5
$ws.Range("H90").Value = 174999.5
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 174999.5
$ws.Range("K90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("M90").Value = 524998.5
$ws.Range("N90").Value = -537478.5
$ws.Range("H92").Value = 1977.2
$ws.Range("I92").Value = 2341.5
$ws.Range("K92").Value = 2341.5
$ws.Range("M92").Value = -1093.5
$ws.Range("H98").Value = 2220.318
$ws.Range("I98").Value = 1662.2354
$ws.Range("J98").Value = 4117.8
$ws.Range("K98").Value = 1662.2354
$ws.Range("L98").Value = 4117.8
$ws.Range("M98").Value = -164.2354
$ws.Range("N98").Value = -7113.8
$ws.Range("H111").Value = 462
$ws.Range("I111").Value = 462
$ws.Range("K111").Value = 1386
$ws.Range("M111").Value = 1681
$ws.Range("H112").Value = 6247.4287
$ws.Range("J112").Value = 6247.4287
$ws.Range("L112").Value = 18742.2861
$ws.Range("N112").Value = -20958.2861
$ws.Range("H116").Value = 6903.857
$ws.Range("I116").Value = 6211.6665
$ws.Range("K116").Value = 6211.6665
$ws.Range("M116").Value = -2769.6665
$ws.Range("H122").Value = 2220.318
$ws.Range("I122").Value = 1662.2354
$ws.Range("J122").Value = 4117.8
$ws.Range("K122").Value = 4986.706200000001
$ws.Range("L122").Value = 12353.4
$ws.Range("M122").Value = -2536.706200000001
$ws.Range("N122").Value = -17253.4
$ws.Range("H132").Value = 2536.9016
$ws.Range("I132").Value = 2394.1035
$ws.Range("J132").Value = 5297.6665
$ws.Range("K132").Value = 7182.310500000001
$ws.Range("L132").Value = 15892.9995
$ws.Range("M132").Value = -4652.310500000001
$ws.Range("N132").Value = -20952.9995
$ws.Range("H137").Value = 2334.373
$ws.Range("I137").Value = 2228.5642
$ws.Range("K137").Value = 6685.692599999999
$ws.Range("M137").Value = -4135.692599999999
$ws.Range("H138").Value = 3587.7778
$ws.Range("I138").Value = 3192.4707
$ws.Range("J138").Value = 3709.9636
$ws.Range("K138").Value = 9577.4121
$ws.Range("L138").Value = 11129.8908
$ws.Range("M138").Value = -4437.4121
$ws.Range("N138").Value = -21409.8908
$ws.Range("H141").Value = 2025.0588
$ws.Range("I141").Value = 2298.2307
$ws.Range("J141").Value = 1137.25
$ws.Range("K141").Value = 6894.6921
$ws.Range("L141").Value = 3411.75
$ws.Range("M141").Value = -1714.6921
$ws.Range("N141").Value = -13771.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2183.3076
$ws.Range("I2").Value = 2353.55
$ws.Range("J2").Value = 1615.8334
$ws.Range("K2").Value = 2353.55
$ws.Range("L2").Value = 1615.8334
$ws.Range("M2").Value = -2240.55
$ws.Range("N2").Value = -1841.8334
$ws.Range("H32").Value = 12615.226
$ws.Range("I32").Value = 8187.5474
$ws.Range("K32").Value = 8187.5474
$ws.Range("M32").Value = -7900.5474
$ws.Range("H45").Value = 6940.732
$ws.Range("I45").Value = 17029.54
$ws.Range("J45").Value = 3890.628
$ws.Range("K45").Value = 17029.54
$ws.Range("L45").Value = 3890.628
$ws.Range("M45").Value = -16652.54
$ws.Range("N45").Value = -4644.628000000001
$ws.Range("H61").Value = 6967.125
$ws.Range("I61").Value = 5105.5
$ws.Range("K61").Value = 5105.5
$ws.Range("M61").Value = -4893.5
$ws.Range("H74").Value = 5192.8076
$ws.Range("I74").Value = 3620.111
$ws.Range("J74").Value = 8731.375
$ws.Range("K74").Value = 3620.111
$ws.Range("L74").Value = 8731.375
$ws.Range("M74").Value = -2746.111
$ws.Range("N74").Value = -10479.375
$ws.Range("H76").Value = 116734.8
$ws.Range("J76").Value = 116734.8
$ws.Range("L76").Value = 116734.8
$ws.Range("N76").Value = -117410.8
$ws.Range("H77").Value = 5192.8076
$ws.Range("I77").Value = 3620.111
$ws.Range("J77").Value = 8731.375
$ws.Range("K77").Value = 18100.555
$ws.Range("L77").Value = 43656.875
$ws.Range("M77").Value = -13732.555
$ws.Range("N77").Value = -52392.875
$ws.Range("H79").Value = 116734.8
$ws.Range("J79").Value = 116734.8
$ws.Range("L79").Value = 116734.8
$ws.Range("N79").Value = -119074.8
$ws.Range("H102").Value = 4900.483
$ws.Range("I102").Value = 5430.8
$ws.Range("J102").Value = 1586
$ws.Range("K102").Value = 5430.8
$ws.Range("L102").Value = 1586
$ws.Range("M102").Value = -3808.8
$ws.Range("N102").Value = -4830
$ws.Range("H110").Value = 5366.8
$ws.Range("I110").Value = 5133.5
$ws.Range("J110").Value = 6300
$ws.Range("K110").Value = 5133.5
$ws.Range("L110").Value = 6300
$ws.Range("M110").Value = -3088.5
$ws.Range("N110").Value = -10390
$ws.Range("H116").Value = 2183.3076
$ws.Range("I116").Value = 2353.55
$ws.Range("J116").Value = 1615.8334
$ws.Range("K116").Value = 2353.55
$ws.Range("L116").Value = 1615.8334
$ws.Range("M116").Value = -59.55000000000018
$ws.Range("N116").Value = -6203.8334
$ws.Range("H122").Value = 10627.521
$ws.Range("I122").Value = 11531.071
$ws.Range("K122").Value = 34593.213
$ws.Range("M122").Value = -32143.213
$ws.Range("H132").Value = 3816.1853
$ws.Range("I132").Value = 2806.6667
$ws.Range("K132").Value = 8420.000100000001
$ws.Range("M132").Value = -5890.000100000001
$ws.Range("H136").Value = 6967.125
$ws.Range("I136").Value = 5105.5
$ws.Range("K136").Value = 15316.5
$ws.Range("M136").Value = -12766.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2183.3076
$ws.Range("I3").Value = 2353.55
$ws.Range("J3").Value = 1615.8334
$ws.Range("K3").Value = 2353.55
$ws.Range("L3").Value = 1615.8334
$ws.Range("M3").Value = -2239.55
$ws.Range("N3").Value = -1843.8334
$ws.Range("H86").Value = 33335500
$ws.Range("I86").Value = 66668268
$ws.Range("J86").Value = 2733
$ws.Range("K86").Value = 66668268
$ws.Range("L86").Value = 2733
$ws.Range("M86").Value = -66667145
$ws.Range("N86").Value = -4979
$ws.Range("H89").Value = 33335500
$ws.Range("I89").Value = 66668268
$ws.Range("J89").Value = 2733
$ws.Range("K89").Value = 333341340
$ws.Range("L89").Value = 13665
$ws.Range("M89").Value = -333335724
$ws.Range("N89").Value = -24897
$ws.Range("H99").Value = 3826.3914
$ws.Range("I99").Value = 3727.5908
$ws.Range("K99").Value = 3727.5908
$ws.Range("M99").Value = -2229.5908
$ws.Range("H107").Value = 2624.25
$ws.Range("I107").Value = 3165.6667
$ws.Range("K107").Value = 3165.6667
$ws.Range("M107").Value = -1245.6667
$ws.Range("H138").Value = 86747
$ws.Range("J138").Value = 86747
$ws.Range("L138").Value = 86747
$ws.Range("N138").Value = -97027

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6119.731
$ws.Range("J31").Value = 5618.1035
$ws.Range("L31").Value = 5618.1035
$ws.Range("N31").Value = -6208.1035
$ws.Range("H34").Value = 6119.731
$ws.Range("J34").Value = 5618.1035
$ws.Range("L34").Value = 5618.1035
$ws.Range("N34").Value = -6022.1035
$ws.Range("H74").Value = 35333
$ws.Range("J74").Value = 35333
$ws.Range("L74").Value = 35333
$ws.Range("N74").Value = -37081
$ws.Range("H77").Value = 35333
$ws.Range("J77").Value = 35333
$ws.Range("L77").Value = 105999
$ws.Range("N77").Value = -114735
$ws.Range("H97").Value = 98197
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 98197
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").Value = 98197
$ws.Range("N97").Value = -100179
$ws.Range("H99").Value = 6102.3335
$ws.Range("I99").Value = 7605
$ws.Range("J99").Value = 4599.6665
$ws.Range("K99").Value = 7605
$ws.Range("L99").Value = 4599.6665
$ws.Range("M99").Value = -6107
$ws.Range("N99").Value = -7595.6665
$ws.Range("H105").Value = 5108.7827
$ws.Range("I105").Value = 2875.15
$ws.Range("K105").Value = 2875.15
$ws.Range("M105").Value = -1128.15
$ws.Range("H107").Value = 2271.6428
$ws.Range("I107").Value = 1960.5
$ws.Range("J107").Value = 3049.5
$ws.Range("K107").Value = 1960.5
$ws.Range("L107").Value = 3049.5
$ws.Range("M107").Value = -40.5
$ws.Range("N107").Value = -6889.5
$ws.Range("H126").Value = 6102.3335
$ws.Range("I126").Value = 7605
$ws.Range("J126").Value = 4599.6665
$ws.Range("K126").Value = 22815
$ws.Range("L126").Value = 13798.9995
$ws.Range("M126").Value = -20345
$ws.Range("N126").Value = -18738.9995
$ws.Range("H132").Value = 2341.5
$ws.Range("I132").Value = 2318.8572
$ws.Range("K132").Value = 6956.571599999999
$ws.Range("M132").Value = -4426.571599999999
$ws.Range("H141").Value = 208075
$ws.Range("J141").Value = 216394.73
$ws.Range("L141").Value = 216394.73
$ws.Range("N141").Value = -226754.73

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 338.625
$ws.Range("J12").Value = 579.25
$ws.Range("L12").Value = 1737.75
$ws.Range("N12").Value = -2083.75
$ws.Range("H34").Value = 4000
$ws.Range("I34").Value = 2250
$ws.Range("J34").Value = 7500
$ws.Range("K34").Value = 6750
$ws.Range("L34").Value = 22500
$ws.Range("M34").Value = -6666
$ws.Range("N34").Value = -22668
$ws.Range("H55").Value = 7604.5
$ws.Range("J55").Value = 9445.111000000001
$ws.Range("L55").Value = 28335.333
$ws.Range("N55").Value = -28689.333
$ws.Range("H113").Value = 3030.6453
$ws.Range("I113").Value = 1679.4
$ws.Range("J113").Value = 3290.5
$ws.Range("K113").Value = 5038.200000000001
$ws.Range("L113").Value = 9871.5
$ws.Range("M113").Value = -2868.200000000001
$ws.Range("N113").Value = -14211.5
$ws.Range("H131").Value = 3254
$ws.Range("I131").Value = 988.8
$ws.Range("J131").Value = 4872
$ws.Range("K131").Value = 2966.4
$ws.Range("L131").Value = 14616
$ws.Range("M131").Value = 2073.6
$ws.Range("N131").Value = -24696
$ws.Range("H141").Value = 8575.764999999999
$ws.Range("I141").Value = 6800.0835
$ws.Range("K141").Value = 20400.2505
$ws.Range("M141").Value = -15220.2505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 14613.286
$ws.Range("J46").Value = 14613.286
$ws.Range("L46").Value = 14613.286
$ws.Range("N46").Value = -14925.286
$ws.Range("H80").Value = 8537.777
$ws.Range("I80").Value = 6012.25
$ws.Range("K80").Value = 6012.25
$ws.Range("M80").Value = -5014.25
$ws.Range("H83").Value = 8537.777
$ws.Range("I83").Value = 6012.25
$ws.Range("K83").Value = 30061.25
$ws.Range("M83").Value = -25069.25
$ws.Range("H97").Value = 875.63635
$ws.Range("I97").Value = 918.2
$ws.Range("J97").Value = 450
$ws.Range("K97").Value = 918.2
$ws.Range("L97").Value = 450
$ws.Range("M97").Value = -422.2
$ws.Range("N97").Value = -1442
$ws.Range("H102").Value = 2678.4546
$ws.Range("I102").Value = 2622.6843
$ws.Range("K102").Value = 2622.6843
$ws.Range("M102").Value = -1000.6843
$ws.Range("H126").Value = 3586.5557
$ws.Range("I126").Value = 3282
$ws.Range("K126").Value = 9846
$ws.Range("M126").Value = -7376
$ws.Range("H132").Value = 4564.755
$ws.Range("I132").Value = 4253.925
$ws.Range("J132").Value = 5946.222
$ws.Range("K132").Value = 12761.775
$ws.Range("L132").Value = 17838.666
$ws.Range("M132").Value = -10231.775
$ws.Range("N132").Value = -22898.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19901
$ws.Range("I40").Value = 29666.666
$ws.Range("J40").Value = 5252.5
$ws.Range("K40").Value = 29666.666
$ws.Range("L40").Value = 5252.5
$ws.Range("M40").Value = -29530.666
$ws.Range("N40").Value = -5524.5
$ws.Range("H46").Value = 937.6667
$ws.Range("I46").Value = 550
$ws.Range("J46").Value = 1131.5
$ws.Range("K46").Value = 550
$ws.Range("L46").Value = 1131.5
$ws.Range("M46").Value = -362
$ws.Range("N46").Value = -1507.5
$ws.Range("H55").Value = 629.1923
$ws.Range("J55").Value = 715.4545000000001
$ws.Range("L55").Value = 715.4545000000001
$ws.Range("N55").Value = -1061.4545
$ws.Range("H100").Value = 3165.1
$ws.Range("I100").Value = 3378.75
$ws.Range("J100").Value = 2310.5
$ws.Range("K100").Value = 3378.75
$ws.Range("L100").Value = 2310.5
$ws.Range("M100").Value = -2837.75
$ws.Range("N100").Value = -3392.5
$ws.Range("H136").Value = 2502.88
$ws.Range("I136").Value = 2230.682
$ws.Range("K136").Value = 6692.045999999999
$ws.Range("M136").Value = -4142.045999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 25937.625
$ws.Range("J69").Value = 24818.572
$ws.Range("L69").Value = 24818.572
$ws.Range("N69").Value = -26316.572
$ws.Range("H72").Value = 25937.625
$ws.Range("J72").Value = 24818.572
$ws.Range("L72").Value = 74455.716
$ws.Range("N72").Value = -81943.716
$ws.Range("H100").Value = 2628.4
$ws.Range("I100").Value = 2901.3
$ws.Range("K100").Value = 5802.6
$ws.Range("M100").Value = -5261.6
$ws.Range("H122").Value = 7497.8
$ws.Range("I122").Value = 6687.5454
$ws.Range("K122").Value = 20062.6362
$ws.Range("M122").Value = -17612.6362
$ws.Range("H123").Value = 84166.414
$ws.Range("J123").Value = 84166.414
$ws.Range("L123").Value = 84166.414
$ws.Range("N123").Value = -93966.414
$ws.Range("H125").Value = 92749.375
$ws.Range("J125").Value = 92749.375
$ws.Range("L125").Value = 92749.375
$ws.Range("N125").Value = -102589.375
$ws.Range("H126").Value = 1947.6522
$ws.Range("I126").Value = 1918
$ws.Range("K126").Value = 5754
$ws.Range("M126").Value = -3284
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").Value = 0
$ws.Range("H132").Value = 3931.575
$ws.Range("I132").Value = 3452.8572
$ws.Range("K132").Value = 10358.5716
$ws.Range("M132").Value = -7828.571599999999
$ws.Range("H136").Value = 15974.171
$ws.Range("I136").Value = 24506.96
$ws.Range("J136").Value = 2641.6875
$ws.Range("K136").Value = 73520.88
$ws.Range("L136").Value = 7925.0625
$ws.Range("M136").Value = -70970.88
$ws.Range("N136").Value = -13025.0625
